# Updates the numeric results table on Sheet1 (rows 2-25, columns B..O)
# with the newly computed values for the "380 kV" case.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns that receive new values, in left-to-right order.
# (Columns G, I, K, M stay at 0 and are left untouched.)
$colLetters = @("B","C","D","E","F","H","J","L","N","O")

# New values per data row (rows 2..25 correspond to A-column index 0..23),
# one array per row, values aligned with $colLetters above.
$newValues = @(
    @(1.606891053220522,0.1429500831398514,0.1249724091649256,0.1209660387166416,1.675643970990237,0.07973214163530429,0.1481101581286079,0.3160357423339661,1.454331322397671,4.329673273218418),
    @(1.510906854829216,0.1320061420179144,0.1243874862339993,0.1215383176541014,1.678269952518939,0.07973214163530429,0.149362862798716,0.3092004679037927,1.467457292950058,4.339556449032727),
    @(1.452408512553063,0.1252392689791151,0.1240625125243611,0.1219274526726313,1.680778051694617,0.07973214163530429,0.1501807542708793,0.3051170192273389,1.476064786199338,4.348080734095419),
    @(1.428681263705585,0.1224699808856116,0.1239387148679256,0.1220955386506208,1.682025394675051,0.07973214163530429,0.1505263150354246,0.3034816613591573,1.479710325444593,4.352171853911216),
    @(1.424748140946463,0.1220094384895276,0.1239186808517232,0.1221240240951449,1.682246121894778,0.07973214163530429,0.1505844361340536,0.3032118482649935,1.480323997060601,4.352888466052889),
    @(1.452088065832015,0.1252019686959045,0.1240608079437635,0.1219296810094086,1.680793961639083,0.07973214163530429,0.1501853649558669,0.3050948478846749,1.476113392675941,4.348133408754222),
    @(1.573706052430623,0.1391865047557275,0.1247636585373897,0.1211555348524183,1.676363514665837,0.07973214163530429,0.1485319842800124,0.313655486449619,1.458743455477588,4.332571142067025),
    @(1.81560421863685,0.1662299329407801,0.1264115142573132,0.1199362711098306,1.674783178183162,0.07973214163530429,0.1456758152411375,0.3313372079275467,1.429025804156133,4.321553114353975),
    @(1.995345918229589,0.1858616385541723,0.1277846476447593,0.1192217142890968,1.67795867846921,0.07973214163530429,0.1438120303164716,0.3448670973745322,1.409835277288416,4.325369774305301),
    @(2.077542605579595,0.1947401106862117,0.1284442136094555,0.1189358017912561,1.680345898222086,0.07973214163530429,0.1430149403843348,0.3511380776938609,1.401677819984691,4.329698269913195),
    @(2.108729037880778,0.1980945399566281,0.1286989594349777,0.118833146769008,1.681385461308736,0.07973214163530429,0.1427203910824453,0.3535293053748632,1.398671050879187,4.331710474109286),
    @(2.102009825291248,0.1973724466967894,0.1286438743862419,0.1188550059240416,1.681155542749977,0.07973214163530429,0.142783503372792,0.3530135784903052,1.399314953254731,4.331260511081894),
    @(2.080107132670037,0.1950162359712238,0.1284650720677547,0.1189272438760582,1.680428706466515,0.07973214163530429,0.14299056154689,0.3513344747983211,1.401428802356286,4.329856335853435),
    @(2.066698907831267,0.1935719870141384,0.1283561981993415,0.1189722223543992,1.680001154694551,0.07973214163530429,0.1431183399961871,0.3503081249893256,1.402734310861781,4.329044836017545),
    @(1.989982705617763,0.1852803472305595,0.1277422428104629,0.1192411857264855,1.677821638668462,0.07973214163530429,0.1438651427259945,0.3444595978459404,1.410379902070765,4.325139088433218),
    @(1.943029061000118,0.180180235888912,0.127374517268521,0.1194162008901056,1.676726030552118,0.07973214163530429,0.1443362764878291,0.340901354057479,1.415216810345854,4.323407218344357),
    @(1.91606326184592,0.1772418936918996,0.1271663004906998,0.1195205501910994,1.67618459295997,0.07973214163530429,0.1446120379054854,0.3388656900117866,1.418052751354878,4.322655070123517),
    @(1.906940146557019,0.1762461861647182,0.127096367987491,0.1195565144640014,1.676016509872767,0.07973214163530429,0.1447062267853045,0.3381783336393909,1.41902220703949,4.322442299007633),
    @(1.948023157047317,0.1807236590567527,0.1274133221408107,0.1193971889301633,1.6768334768242,0.07973214163530429,0.1442856290688219,0.3412790037225335,1.414696337119743,4.323566325699545),
    @(2.086538860493818,0.1957085212566199,0.1285174557291953,0.1189058735988304,1.680638516279686,0.07973214163530429,0.1429295457132493,0.3518272205487989,1.400805680704742,4.330258647778152),
    @(2.177417753308873,0.2054572969294952,0.1292680947175313,0.11861748550794,1.683915572558959,0.07973214163530429,0.1420857632085628,0.3588174439680643,1.392206931780059,4.336807314771505),
    @(2.12888242950271,0.2002583361715438,0.1288648217406489,0.1187684151412896,1.682094231954196,0.07973214163530429,0.1425322197324874,0.3550778668584798,1.396752370912438,4.333113059959146),
    @(1.945765237391299,0.180477996806701,0.1273957685084497,0.1194057726154192,1.676784624871615,0.07973214163530429,0.1443085115019311,0.3411082370938061,1.414931471307661,4.323493634706921),
    @(1.74980567256415,0.158955224198877,0.1259370616267645,0.1202342195758757,1.674449407329192,0.07973214163530429,0.1464072261233316,0.3264587194022539,1.436600684554776,4.322443730337682)
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowNum = $startRow + $i
    $rowVals = $newValues[$i]
    for ($j = 0; $j -lt $colLetters.Count; $j++) {
        $cellRef = "$($colLetters[$j])$rowNum"
        $ws.Range($cellRef).Value = $rowVals[$j]
    }
}
